# Ajout des feedbacks "board"
# On the "Feedback" worksheet, mark rows 8, 10, 11 and 12 in column H
# (the "Done" column) with an "X", matching the existing rows that
# already have this feedback flag, and update the active selection to H12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feedback")

$ws.Range("H8").Value = "X"
$ws.Range("H10").Value = "X"
$ws.Range("H11").Value = "X"
$ws.Range("H12").Value = "X"

$ws.Range("H12").Select() | Out-Null
